# Updates the cryptos list table (columns B:E, rows 2-51) to the latest
# scraped values, matching the GitHub Actions data refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '67.856.87'
$ws.Range("E2").Value = '  -7.35%  '

# Row 3
$ws.Range("D3").Value = '3.680.28'
$ws.Range("E3").Value = '  -7.65%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.05%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '572.28'
$ws.Range("E5").Value = '  -6.21%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.06'
$ws.Range("E6").Value = '  -0.68%  '

# Row 7
$ws.Range("D7").Value = '3.665.32'
$ws.Range("E7").Value = '  -7.78%  '

# Row 8
$ws.Range("E8").Value = '  -10.76%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.01'
$ws.Range("E9").Value = '  +0.54%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.699'
$ws.Range("E10").Value = '  -12.78%  '

# Row 11
$ws.Range("E11").Value = '  -13.39%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '51.05'
$ws.Range("E12").Value = '  -10.15%  '

# Row 13
$ws.Range("E13").Value = '  -13.82%  '

# Row 14
$ws.Range("E14").Value = '  -11.35%  '

# Row 15
$ws.Range("D15").Value = '4.265.84'
$ws.Range("E15").Value = '  -7.69%  '

# Row 16
$ws.Range("D16").Value = '3.676.37'
$ws.Range("E16").Value = '  -7.83%  '

# Row 17
$ws.Range("B17").Value = 'Chainlink'
$ws.Range("C17").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '19.30'
$ws.Range("E17").Value = '  -9.02%  '

# Row 18
$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.126'
$ws.Range("E18").Value = '  -3.48%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.80'
$ws.Range("E19").Value = '  -10.73%  '

# Row 20
$ws.Range("E20").Value = '  -11.16%  '

# Row 21
$ws.Range("D21").Value = '67.696.74'
$ws.Range("E21").Value = '  -7.50%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '404.76'
$ws.Range("E22").Value = '  -13.07%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.45'
$ws.Range("E23").Value = '  -8.09%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '87.34'
$ws.Range("E24").Value = '  -10.67%  '

# Row 25
$ws.Range("E25").Value = '  -11.39%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.71'
$ws.Range("E26").Value = '  -11.35%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.57'
$ws.Range("E27").Value = '  -7.23%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.98'
$ws.Range("E28").Value = '  +1.79%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.72'
$ws.Range("E29").Value = '  -12.54%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.41'

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '32.54'
$ws.Range("E31").Value = '  -10.98%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.51'
$ws.Range("E32").Value = '  -5.68%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '12.43'
$ws.Range("E33").Value = '  -12.64%  '

# Row 34
$ws.Range("E34").Value = '  -10.99%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '64.87'
$ws.Range("E35").Value = '  -8.23%  '

# Row 36
$ws.Range("B36").Value = 'Bittensor'
$ws.Range("C36").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '597.07'
$ws.Range("E36").Value = '  -7.14%  '

# Row 37
$ws.Range("B37").Value = 'InjectiveProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '42.81'
$ws.Range("E37").Value = '  -14.19%  '

# Row 38
$ws.Range("D38").Value = '0.0₃0892'
$ws.Range("E38").Value = '  -13.75%  '

# Row 39
$ws.Range("E39").Value = '  +0.05%  '

# Row 40
$ws.Range("E40").Value = '  -9.56%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.999'
$ws.Range("E41").Value = '  -0.15%  '

# Row 42
$ws.Range("E42").Value = '  -10.22%  '

# Row 43
$ws.Range("B43").Value = 'ThetaToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.96'
$ws.Range("E43").Value = '  -13.69%  '

# Row 44
$ws.Range("B44").Value = 'dogwifhat'
$ws.Range("C44").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.94'
$ws.Range("E44").Value = '  -13.21%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0435'
$ws.Range("E45").Value = '  -10.79%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.60'
$ws.Range("E46").Value = '  -3.25%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.13'
$ws.Range("E47").Value = '  -13.68%  '

# Row 48
$ws.Range("B48").Value = 'WEMIXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.71'
$ws.Range("E48").Value = '  -10.73%  '

# Row 49
$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.133'
$ws.Range("E49").Value = '  -11.92%  '

# Row 50
$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").Value = '2.732.60'
$ws.Range("E50").Value = '  -3.46%  '

# Row 51
$ws.Range("B51").Value = 'ApeXProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.11'
$ws.Range("E51").Value = '  -9.33%  '

